$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.380.26"
$ws.Range("E2").Value = "  -3.41%  "

$ws.Range("D3").Value = "1.747.71"
$ws.Range("E3").Value = "  -3.64%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.91"
$ws.Range("E5").Value = "  -2.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4238"
$ws.Range("E7").Value = "  -4.85%  "

$ws.Range("E8").Value = "  -3.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07489"
$ws.Range("E9").Value = "  -3.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.13"
$ws.Range("E10").Value = "  -6.32%  "

$ws.Range("E11").Value = "  -2.87%  "

$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.66"
$ws.Range("E13").Value = "  -6.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.027"
$ws.Range("E14").Value = "  -4.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.216"
$ws.Range("E15").Value = "  -4.91%  "

$ws.Range("D16").Value = "1.743.48"
$ws.Range("E16").Value = "  -5.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.44"
$ws.Range("E17").Value = "  -0.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001066"
$ws.Range("E18").Value = "  -1.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06368"
$ws.Range("E19").Value = "  -2.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.06"
$ws.Range("E21").Value = "  -2.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.895"
$ws.Range("E22").Value = "  -5.64%  "

$ws.Range("D23").Value = "27.439.39"
$ws.Range("E23").Value = "  -3.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  -4.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.094"
$ws.Range("E25").Value = "  -2.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.89"
$ws.Range("E26").Value = "  +3.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.23"
$ws.Range("E27").Value = "  -2.74%  "

$ws.Range("D28").Value = "1.948.14"
$ws.Range("E28").Value = "  -4.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.117"
$ws.Range("E29").Value = "  -8.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.80"
$ws.Range("E30").Value = "  -3.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.101"
$ws.Range("E31").Value = "  -8.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.650"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.546"
$ws.Range("E33").Value = "  -6.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08902"
$ws.Range("E34").Value = "  -3.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.19"
$ws.Range("E35").Value = "  -6.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02281"
$ws.Range("E36").Value = "  -3.21%  "

$ws.Range("E37").Value = "  -4.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05994"
$ws.Range("E38").Value = "  -3.83%  "

$ws.Range("E39").Value = "  -3.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.949"
$ws.Range("E40").Value = "  -4.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.181"
$ws.Range("E41").Value = "  -1.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9996"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.892"
$ws.Range("E43").Value = "  -3.17%  "

$ws.Range("E44").Value = "  -1.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.30"
$ws.Range("E45").Value = "  -4.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5882"
$ws.Range("E46").Value = "  -3.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.688"
$ws.Range("E47").Value = "  -2.14%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.58"
$ws.Range("E48").Value = "  -3.60%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.969"
$ws.Range("E49").Value = "  -3.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.163"
$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06809"
$ws.Range("E51").Value = "  -2.57%  "
